# Applies the "pembaruan format interpretasi hasil" edit:
#  1) Row 51 (age + q1..q42 answer columns) switches from text to real numbers.
#  2) A new response row (52) is appended, keeping its values as text,
#     the same way the earlier rows were originally recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 51: reinterpret age (col B) and q1..q42 (col G..AV) as numbers ---
$ws.Cells.Item(51, 2).Value = 34
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 2
$ws.Cells.Item(51, 9).Value = 3
$ws.Cells.Item(51, 10).Value = 2
$ws.Cells.Item(51, 11).Value = 3
$ws.Cells.Item(51, 12).Value = 2
$ws.Cells.Item(51, 13).Value = 2
$ws.Cells.Item(51, 14).Value = 1
$ws.Cells.Item(51, 15).Value = 2
$ws.Cells.Item(51, 16).Value = 2
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = 2
$ws.Cells.Item(51, 19).Value = 1
$ws.Cells.Item(51, 20).Value = 2
$ws.Cells.Item(51, 21).Value = 1
$ws.Cells.Item(51, 22).Value = 2
$ws.Cells.Item(51, 23).Value = 3
$ws.Cells.Item(51, 24).Value = 2
$ws.Cells.Item(51, 25).Value = 2
$ws.Cells.Item(51, 26).Value = 1
$ws.Cells.Item(51, 27).Value = 2
$ws.Cells.Item(51, 28).Value = 3
$ws.Cells.Item(51, 29).Value = 2
$ws.Cells.Item(51, 30).Value = 2
$ws.Cells.Item(51, 31).Value = 1
$ws.Cells.Item(51, 32).Value = 2
$ws.Cells.Item(51, 33).Value = 2
$ws.Cells.Item(51, 34).Value = 1
$ws.Cells.Item(51, 35).Value = 2
$ws.Cells.Item(51, 36).Value = 2
$ws.Cells.Item(51, 37).Value = 1
$ws.Cells.Item(51, 38).Value = 2
$ws.Cells.Item(51, 39).Value = 1
$ws.Cells.Item(51, 40).Value = 2
$ws.Cells.Item(51, 41).Value = 3
$ws.Cells.Item(51, 42).Value = 2
$ws.Cells.Item(51, 43).Value = 2
$ws.Cells.Item(51, 44).Value = 1
$ws.Cells.Item(51, 45).Value = 2
$ws.Cells.Item(51, 46).Value = 2
$ws.Cells.Item(51, 47).Value = 3
$ws.Cells.Item(51, 48).Value = 2

# --- 2) Row 52: new response, entered as text (matches how prior rows were stored) ---
$ws.Range("B52").NumberFormat = "@"
$ws.Range("G52:AV52").NumberFormat = "@"

$ws.Cells.Item(52, 1).Value = '2025-05-21 10:34:00'
$ws.Cells.Item(52, 2).Value = '34'
$ws.Cells.Item(52, 3).Value = 'Jakarta'
$ws.Cells.Item(52, 4).Value = 'SMA/SMK'
$ws.Cells.Item(52, 5).Value = 'male'
$ws.Cells.Item(52, 6).Value = 'dasdsd'
$ws.Cells.Item(52, 7).Value = '2'
$ws.Cells.Item(52, 8).Value = '3'
$ws.Cells.Item(52, 9).Value = '2'
$ws.Cells.Item(52, 10).Value = '2'
$ws.Cells.Item(52, 11).Value = '3'
$ws.Cells.Item(52, 12).Value = '2'
$ws.Cells.Item(52, 13).Value = '2'
$ws.Cells.Item(52, 14).Value = '2'
$ws.Cells.Item(52, 15).Value = '1'
$ws.Cells.Item(52, 16).Value = '2'
$ws.Cells.Item(52, 17).Value = '3'
$ws.Cells.Item(52, 18).Value = '2'
$ws.Cells.Item(52, 19).Value = '2'
$ws.Cells.Item(52, 20).Value = '3'
$ws.Cells.Item(52, 21).Value = '2'
$ws.Cells.Item(52, 22).Value = '1'
$ws.Cells.Item(52, 23).Value = '2'
$ws.Cells.Item(52, 24).Value = '3'
$ws.Cells.Item(52, 25).Value = '2'
$ws.Cells.Item(52, 26).Value = '3'
$ws.Cells.Item(52, 27).Value = '2'
$ws.Cells.Item(52, 28).Value = '2'
$ws.Cells.Item(52, 29).Value = '3'
$ws.Cells.Item(52, 30).Value = '2'
$ws.Cells.Item(52, 31).Value = '1'
$ws.Cells.Item(52, 32).Value = '2'
$ws.Cells.Item(52, 33).Value = '3'
$ws.Cells.Item(52, 34).Value = '2'
$ws.Cells.Item(52, 35).Value = '2'
$ws.Cells.Item(52, 36).Value = '2'
$ws.Cells.Item(52, 37).Value = '3'
$ws.Cells.Item(52, 38).Value = '2'
$ws.Cells.Item(52, 39).Value = '2'
$ws.Cells.Item(52, 40).Value = '2'
$ws.Cells.Item(52, 41).Value = '1'
$ws.Cells.Item(52, 42).Value = '2'
$ws.Cells.Item(52, 43).Value = '2'
$ws.Cells.Item(52, 44).Value = '3'
$ws.Cells.Item(52, 45).Value = '2'
$ws.Cells.Item(52, 46).Value = '1'
$ws.Cells.Item(52, 47).Value = '2'
$ws.Cells.Item(52, 48).Value = '2'
